# Auto-generated edit script: updates Leve market-price/profit data cells
# across multiple sheets to reflect a scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3319.087
$ws.Range("I15").Value = 3319.087
$ws.Range("K15").Value = 9957.261
$ws.Range("M15").Value = -9788.261
$ws.Range("H17").Value = 685.7727
$ws.Range("J17").Value = 671.87805
$ws.Range("L17").Value = 2015.63415
$ws.Range("N17").Value = -2351.63415
$ws.Range("H58").Value = 324.07693
$ws.Range("J58").Value = 900
$ws.Range("L58").Value = 2700
$ws.Range("N58").Value = -3000
$ws.Range("H94").Value = 10374.875
$ws.Range("I94").Value = 4714.143
$ws.Range("J94").Value = 50000
$ws.Range("K94").Value = 4714.143
$ws.Range("L94").Value = 50000
$ws.Range("M94").Value = -4263.143
$ws.Range("N94").Value = -50902
$ws.Range("H125").Value = 4241
$ws.Range("I125").Value = 5250.5
$ws.Range("K125").Value = 47254.5
$ws.Range("M125").Value = -44794.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 80937.28999999999
$ws.Range("I32").Value = 80937.28999999999
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 80937.28999999999
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -80650.28999999999
$ws.Range("N32").ClearContents()
$ws.Range("H37").Value = 66916.586
$ws.Range("J37").Value = 78499.875
$ws.Range("L37").Value = 78499.875
$ws.Range("N37").Value = -79045.875
$ws.Range("H61").Value = 33333332
$ws.Range("I61").Value = 33333332
$ws.Range("K61").Value = 33333332
$ws.Range("M61").Value = -33333120
$ws.Range("H74").Value = 4516.8667
$ws.Range("I74").Value = 1255.2858
$ws.Range("J74").Value = 7370.75
$ws.Range("K74").Value = 1255.2858
$ws.Range("L74").Value = 7370.75
$ws.Range("M74").Value = -381.2858000000001
$ws.Range("N74").Value = -9118.75
$ws.Range("H77").Value = 4516.8667
$ws.Range("I77").Value = 1255.2858
$ws.Range("J77").Value = 7370.75
$ws.Range("K77").Value = 6276.429
$ws.Range("L77").Value = 36853.75
$ws.Range("M77").Value = -1908.429
$ws.Range("N77").Value = -45589.75
$ws.Range("H88").Value = 940.3333
$ws.Range("I88").Value = 800
$ws.Range("J88").Value = 968.4
$ws.Range("K88").Value = 800
$ws.Range("L88").Value = 968.4
$ws.Range("M88").Value = -394
$ws.Range("N88").Value = -1780.4
$ws.Range("H91").Value = 940.3333
$ws.Range("I91").Value = 800
$ws.Range("J91").Value = 968.4
$ws.Range("K91").Value = 800
$ws.Range("L91").Value = 968.4
$ws.Range("M91").Value = 604
$ws.Range("N91").Value = -3776.4
$ws.Range("H136").Value = 33333332
$ws.Range("I136").Value = 33333332
$ws.Range("K136").Value = 99999996
$ws.Range("M136").Value = -99997446

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H48").Value = 400589.5
$ws.Range("J48").Value = 400589.5
$ws.Range("L48").Value = 400589.5
$ws.Range("N48").Value = -401419.5
$ws.Range("H86").Value = 1161.3823
$ws.Range("I86").Value = 1176.0303
$ws.Range("K86").Value = 1176.0303
$ws.Range("M86").Value = -53.0302999999999
$ws.Range("H89").Value = 1161.3823
$ws.Range("I89").Value = 1176.0303
$ws.Range("K89").Value = 5880.1515
$ws.Range("M89").Value = -264.1514999999999
$ws.Range("H94").Value = 1410.2
$ws.Range("I94").Value = 1127.4286
$ws.Range("K94").Value = 1127.4286
$ws.Range("M94").Value = -676.4286
$ws.Range("H99").Value = 3495.8
$ws.Range("I99").Value = 3876.2307
$ws.Range("K99").Value = 3876.2307
$ws.Range("M99").Value = -2378.2307
$ws.Range("H125").Value = 78000
$ws.Range("J125").Value = 78000
$ws.Range("L125").Value = 78000
$ws.Range("N125").Value = -87840
$ws.Range("H134").Value = 694758.9
$ws.Range("I134").Value = 746339.4399999999
$ws.Range("K134").Value = 2239018.32
$ws.Range("M134").Value = -2236483.32

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 419.5
$ws.Range("I7").Value = 424.42856
$ws.Range("K7").Value = 424.42856
$ws.Range("M7").Value = -311.42856
$ws.Range("H58").Value = 1129714.8
$ws.Range("J58").Value = 12198.8
$ws.Range("L58").Value = 12198.8
$ws.Range("N58").Value = -12604.8
$ws.Range("H86").Value = 83046.5
$ws.Range("I86").Value = 5811.077
$ws.Range("J86").Value = 160281.92
$ws.Range("K86").Value = 5811.077
$ws.Range("L86").Value = 160281.92
$ws.Range("M86").Value = -4688.077
$ws.Range("N86").Value = -162527.92
$ws.Range("H89").Value = 83046.5
$ws.Range("I89").Value = 5811.077
$ws.Range("J89").Value = 160281.92
$ws.Range("K89").Value = 29055.385
$ws.Range("L89").Value = 801409.6000000001
$ws.Range("M89").Value = -23439.385
$ws.Range("N89").Value = -812641.6000000001
$ws.Range("H107").Value = 980.2632
$ws.Range("I107").Value = 851.625
$ws.Range("J107").Value = 1666.3334
$ws.Range("K107").Value = 851.625
$ws.Range("L107").Value = 1666.3334
$ws.Range("M107").Value = 1068.375
$ws.Range("N107").Value = -5506.3334
$ws.Range("H136").Value = 1129714.8
$ws.Range("J136").Value = 12198.8
$ws.Range("L136").Value = 36596.39999999999
$ws.Range("N136").Value = -41696.39999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 421.625
$ws.Range("I5").Value = 254.8
$ws.Range("J5").Value = 699.6667
$ws.Range("K5").Value = 764.4000000000001
$ws.Range("L5").Value = 2099.0001
$ws.Range("M5").Value = -652.4000000000001
$ws.Range("N5").Value = -2323.0001
$ws.Range("H22").Value = 734.2
$ws.Range("I22").Value = 150
$ws.Range("J22").Value = 1610.5
$ws.Range("K22").Value = 450
$ws.Range("L22").Value = 4831.5
$ws.Range("M22").Value = -281
$ws.Range("N22").Value = -5169.5
$ws.Range("H24").Value = 685.875
$ws.Range("I24").Value = 533
$ws.Range("J24").Value = 777.6
$ws.Range("K24").Value = 1599
$ws.Range("L24").Value = 2332.8
$ws.Range("M24").Value = -1369
$ws.Range("N24").Value = -2792.8
$ws.Range("H27").Value = 734.2
$ws.Range("I27").Value = 150
$ws.Range("J27").Value = 1610.5
$ws.Range("K27").Value = 450
$ws.Range("L27").Value = 4831.5
$ws.Range("M27").Value = -348
$ws.Range("N27").Value = -5035.5
$ws.Range("H113").Value = 1211.8223
$ws.Range("I113").Value = 333
$ws.Range("J113").Value = 1321.675
$ws.Range("K113").Value = 999
$ws.Range("L113").Value = 3965.025
$ws.Range("M113").Value = 1171
$ws.Range("N113").Value = -8305.025
$ws.Range("H135").Value = 421.625
$ws.Range("I135").Value = 254.8
$ws.Range("J135").Value = 699.6667
$ws.Range("K135").Value = 2293.2
$ws.Range("L135").Value = 6297.0003
$ws.Range("M135").Value = 241.7999999999997
$ws.Range("N135").Value = -11367.0003

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 33000
$ws.Range("J49").Value = 33000
$ws.Range("L49").Value = 33000
$ws.Range("N49").Value = -33368
$ws.Range("H55").Value = 29333.334
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H122").Value = 34797.06
$ws.Range("I122").Value = 53046
$ws.Range("K122").Value = 159138
$ws.Range("M122").Value = -156688
$ws.Range("H126").Value = 1114330.4
$ws.Range("I126").Value = 1854703.1
$ws.Range("J126").Value = 3771.3333
$ws.Range("K126").Value = 5564109.300000001
$ws.Range("L126").Value = 11313.9999
$ws.Range("M126").Value = -5561639.300000001
$ws.Range("N126").Value = -16253.9999
$ws.Range("H132").Value = 40488068
$ws.Range("I132").Value = 59536964
$ws.Range("J132").Value = 9162
$ws.Range("K132").Value = 178610892
$ws.Range("L132").Value = 27486
$ws.Range("M132").Value = -178608362
$ws.Range("N132").Value = -32546

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2181914
$ws.Range("I132").Value = 2905968.8
$ws.Range("K132").Value = 8717906.399999999
$ws.Range("M132").Value = -8715376.399999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3774.9285
$ws.Range("I81").Value = 3350
$ws.Range("K81").Value = 6700
$ws.Range("M81").Value = -5639
$ws.Range("H84").Value = 3774.9285
$ws.Range("I84").Value = 3350
$ws.Range("K84").Value = 33500
$ws.Range("M84").Value = -28196
$ws.Range("H122").Value = 2778.4333
$ws.Range("I122").Value = 2340.125
$ws.Range("K122").Value = 7020.375
$ws.Range("M122").Value = -4570.375
$ws.Range("H126").Value = 4553.9546
$ws.Range("I126").Value = 4352.0557
$ws.Range("K126").Value = 13056.1671
$ws.Range("M126").Value = -10586.1671
$ws.Range("H136").Value = 13102.904
$ws.Range("I136").Value = 13107.925
$ws.Range("K136").Value = 39323.77499999999
$ws.Range("M136").Value = -36773.77499999999

